# Automatische test-sync: 2025-06-26 23:13:50
#
# Adds the new "Testmail #2" row (row 34) to the Logs sheet, extends the
# conditional-formatting ranges that covered the data block, and bumps the
# "Productinformatie" tally on the Dashboard sheet from 2 to 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Append the new log entry on row 34 -----------------------------------
$ws.Range("A34").Value = "Wat zijn de verzendkosten?"
$ws.Range("B34").Value = "mailmind.test@zohomail.eu"
$ws.Range("C34").Value = "Testmail #2: Wat zijn de verzendkosten?"
$ws.Range("D34").Value = "Productinformatie"
$ws.Range("E34").Value = "Beste afzender,`nDank u voor uw interesse in onze producten/diensten. Om u nauwkeurige informatie te verstrekken over de verzendkosten, hebben we meer details nodig zoals het product/dienst waar u naar informeert en het afleveradres. Zou u ons kunnen voorzien van deze informatie zodat we u een precieze schatting van de verzendkosten kunnen geven?`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$ws.Range("F34").Value = "2025-06-26 23:13:07"
$ws.Range("G34").Value = "Ja"
$ws.Range("H34").Value = "Nee"
$ws.Range("I34").Value = "Ja"

# Keep the row at the sheet's standard height (writing the multi-line
# answer into E34 otherwise leaves the row flagged with a custom height).
$ws.Rows.Item(34).AutoFit()

# --- Extend the conditional formatting ranges to include row 34 -----------
$ws.Range("D2:D34").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D34"))
$ws.Range("G2:G34").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G34"))
$ws.Range("H2:H34").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H34"))
$ws.Range("I2:I34").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I34"))

# --- Update the Dashboard "Productinformatie" count ------------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B4").Value = 3
